$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rows = @(7, 30, 52, 68, 98, 116, 148, 162, 186, 206, 224, 247)

foreach ($r in $rows) {
    $ws.Range("B$r").Value = "Dulce de Leche"
}
